$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H11").Value = 4682
$ws_ALC.Range("I11").Value = 4682
$ws_ALC.Range("K11").Value = 4682
$ws_ALC.Range("M11").Value = -4542
$ws_ALC.Range("H113").Value = 2833.3333
$ws_ALC.Range("I113").Value = 2833.3333
$ws_ALC.Range("K113").Value = 2833.3333
$ws_ALC.Range("M113").Value = 420.6667000000002
$ws_ALC.Range("H135").Value = 7754.6772
$ws_ALC.Range("I135").Value = 5299.9
$ws_ALC.Range("J135").Value = 12217.909
$ws_ALC.Range("K135").Value = 47699.1
$ws_ALC.Range("L135").Value = 109961.181
$ws_ALC.Range("M135").Value = -45164.1
$ws_ALC.Range("N135").Value = -115031.181
$ws_ALC.Range("H137").Value = 8305.691999999999
$ws_ALC.Range("J137").Value = 7566.3
$ws_ALC.Range("L137").Value = 22698.9
$ws_ALC.Range("N137").Value = -27798.9
$ws_ALC.Range("H138").Value = 4644.0938
$ws_ALC.Range("I138").Value = 2362.4285
$ws_ALC.Range("K138").Value = 7087.2855
$ws_ALC.Range("M138").Value = -1947.2855
$ws_ALC.Range("H141").Value = 2061.7144
$ws_ALC.Range("I141").Value = 2061.7144
$ws_ALC.Range("J141").Value = 0
$ws_ALC.Range("K141").Value = 6185.1432
$ws_ALC.Range("L141").Value = 0
$ws_ALC.Range("M141").Value = -1005.1432
$ws_ALC.Range("N141").ClearContents()
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 2139395.5
$ws_ARM.Range("I32").Value = 2176.0625
$ws_ARM.Range("J32").Value = 27786028
$ws_ARM.Range("K32").Value = 2176.0625
$ws_ARM.Range("L32").Value = 27786028
$ws_ARM.Range("M32").Value = -1889.0625
$ws_ARM.Range("N32").Value = -27786602
$ws_ARM.Range("H41").Value = 14962
$ws_ARM.Range("I41").Value = 0
$ws_ARM.Range("J41").Value = 14962
$ws_ARM.Range("K41").Value = 0
$ws_ARM.Range("L41").Value = 14962
$ws_ARM.Range("M41").ClearContents()
$ws_ARM.Range("N41").Value = -15790
$ws_ARM.Range("H42").Value = 28000
$ws_ARM.Range("J42").Value = 31000
$ws_ARM.Range("L42").Value = 31000
$ws_ARM.Range("N42").Value = -31972
$ws_ARM.Range("H132").Value = 1017578.3
$ws_ARM.Range("I132").Value = 1234747.9
$ws_ARM.Range("J132").Value = 148900
$ws_ARM.Range("K132").Value = 3704243.7
$ws_ARM.Range("L132").Value = 446700
$ws_ARM.Range("M132").Value = -3701713.7
$ws_ARM.Range("N132").Value = -451760
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 142857140
$ws_BSM.Range("I20").Value = 142857140
$ws_BSM.Range("J20").Value = 0
$ws_BSM.Range("K20").Value = 142857140
$ws_BSM.Range("L20").Value = 0
$ws_BSM.Range("M20").Value = -142856893
$ws_BSM.Range("N20").ClearContents()
$ws_BSM.Range("H81").Value = 16852.666
$ws_BSM.Range("J81").Value = 16852.666
$ws_BSM.Range("L81").Value = 16852.666
$ws_BSM.Range("N81").Value = -18974.666
$ws_BSM.Range("H84").Value = 16852.666
$ws_BSM.Range("J84").Value = 16852.666
$ws_BSM.Range("L84").Value = 50557.99800000001
$ws_BSM.Range("N84").Value = -61165.99800000001
$ws_BSM.Range("H107").Value = 3854650.2
$ws_BSM.Range("J107").Value = 11452.3
$ws_BSM.Range("L107").Value = 11452.3
$ws_BSM.Range("N107").Value = -15292.3
$ws_BSM.Range("H134").Value = 1115534.6
$ws_BSM.Range("I134").Value = 1319198.9
$ws_BSM.Range("J134").Value = 9928.714
$ws_BSM.Range("K134").Value = 3957596.7
$ws_BSM.Range("L134").Value = 29786.142
$ws_BSM.Range("M134").Value = -3955061.7
$ws_BSM.Range("N134").Value = -34856.142
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H41").Value = 15050
$ws_CRP.Range("I41").Value = 100
$ws_CRP.Range("J41").Value = 30000
$ws_CRP.Range("K41").Value = 100
$ws_CRP.Range("L41").Value = 30000
$ws_CRP.Range("M41").Value = 328
$ws_CRP.Range("N41").Value = -30856
$ws_CRP.Range("H42").Value = 34763
$ws_CRP.Range("J42").Value = 34763
$ws_CRP.Range("L42").Value = 34763
$ws_CRP.Range("N42").Value = -35949
$ws_CRP.Range("H141").Value = 735000
$ws_CRP.Range("J141").Value = 713333.3
$ws_CRP.Range("L141").Value = 713333.3
$ws_CRP.Range("N141").Value = -723693.3
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H64").Value = 14536.308
$ws_CUL.Range("J64").Value = 14588.637
$ws_CUL.Range("L64").Value = 43765.911
$ws_CUL.Range("N64").Value = -44305.911
$ws_CUL.Range("H67").Value = 14536.308
$ws_CUL.Range("J67").Value = 14588.637
$ws_CUL.Range("L67").Value = 43765.911
$ws_CUL.Range("N67").Value = -45637.911
$ws_CUL.Range("H113").Value = 722.1667
$ws_CUL.Range("J113").Value = 878.75
$ws_CUL.Range("L113").Value = 2636.25
$ws_CUL.Range("N113").Value = -6976.25
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H41").Value = 6820.3335
$ws_GSM.Range("I41").Value = 16333
$ws_GSM.Range("J41").Value = 3649.4443
$ws_GSM.Range("K41").Value = 16333
$ws_GSM.Range("L41").Value = 3649.4443
$ws_GSM.Range("M41").Value = -15978
$ws_GSM.Range("N41").Value = -4359.4443
$ws_GSM.Range("H43").Value = 11954.444
$ws_GSM.Range("I43").Value = 11954.444
$ws_GSM.Range("K43").Value = 11954.444
$ws_GSM.Range("M43").Value = -11803.444
$ws_GSM.Range("H70").Value = 9675.1
$ws_GSM.Range("J70").Value = 9399
$ws_GSM.Range("L70").Value = 9399
$ws_GSM.Range("N70").Value = -9939
$ws_GSM.Range("H73").Value = 9675.1
$ws_GSM.Range("J73").Value = 9399
$ws_GSM.Range("L73").Value = 9399
$ws_GSM.Range("N73").Value = -11271
$ws_GSM.Range("H102").Value = 923361.3
$ws_GSM.Range("I102").Value = 1455205.8
$ws_GSM.Range("J102").Value = 7406.9443
$ws_GSM.Range("K102").Value = 1455205.8
$ws_GSM.Range("L102").Value = 7406.9443
$ws_GSM.Range("M102").Value = -1453583.8
$ws_GSM.Range("N102").Value = -10650.9443
$ws_GSM.Range("H132").Value = 29415282
$ws_GSM.Range("I132").Value = 45457852
$ws_GSM.Range("J132").Value = 3909.1667
$ws_GSM.Range("K132").Value = 136373556
$ws_GSM.Range("L132").Value = 11727.5001
$ws_GSM.Range("M132").Value = -136371026
$ws_GSM.Range("N132").Value = -16787.5001
$ws_GSM.Range("H133").Value = 75666.664
$ws_GSM.Range("J133").Value = 75666.664
$ws_GSM.Range("L133").Value = 75666.664
$ws_GSM.Range("N133").Value = -85786.664
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H41").Value = 23163
$ws_LTW.Range("I41").Value = 49499
$ws_LTW.Range("J41").Value = 9995
$ws_LTW.Range("K41").Value = 49499
$ws_LTW.Range("L41").Value = 9995
$ws_LTW.Range("M41").Value = -49061
$ws_LTW.Range("N41").Value = -10871
$ws_LTW.Range("H42").Value = 19623.125
$ws_LTW.Range("J42").Value = 19712.857
$ws_LTW.Range("L42").Value = 19712.857
$ws_LTW.Range("N42").Value = -20838.857
$ws_LTW.Range("H43").Value = 6924
$ws_LTW.Range("I43").Value = 23670.666
$ws_LTW.Range("J43").Value = 1900
$ws_LTW.Range("K43").Value = 23670.666
$ws_LTW.Range("L43").Value = 1900
$ws_LTW.Range("M43").Value = -23477.666
$ws_LTW.Range("N43").Value = -2286
$ws_LTW.Range("H46").Value = 31252058
$ws_LTW.Range("I46").Value = 1144
$ws_LTW.Range("K46").Value = 1144
$ws_LTW.Range("M46").Value = -956
$ws_LTW.Range("H49").Value = 19623.125
$ws_LTW.Range("J49").Value = 19712.857
$ws_LTW.Range("L49").Value = 19712.857
$ws_LTW.Range("N49").Value = -20006.857
$ws_LTW.Range("H61").Value = 7414
$ws_LTW.Range("I61").Value = 6615.2144
$ws_LTW.Range("J61").Value = 8430.637000000001
$ws_LTW.Range("K61").Value = 6615.2144
$ws_LTW.Range("L61").Value = 8430.637000000001
$ws_LTW.Range("M61").Value = -6413.2144
$ws_LTW.Range("N61").Value = -8834.637000000001
$ws_LTW.Range("H101").Value = 32761.8
$ws_LTW.Range("J101").Value = 32761.8
$ws_LTW.Range("L101").Value = 32761.8
$ws_LTW.Range("N101").Value = -39251.8
$ws_LTW.Range("H113").Value = 7414
$ws_LTW.Range("I113").Value = 6615.2144
$ws_LTW.Range("J113").Value = 8430.637000000001
$ws_LTW.Range("K113").Value = 6615.2144
$ws_LTW.Range("L113").Value = 8430.637000000001
$ws_LTW.Range("M113").Value = -4445.2144
$ws_LTW.Range("N113").Value = -12770.637
$ws_LTW.Range("H121").Value = 72222
$ws_LTW.Range("J121").Value = 72222
$ws_LTW.Range("L121").Value = 72222
$ws_LTW.Range("N121").Value = -75716
$ws_LTW.Range("H122").Value = 5027.1875
$ws_LTW.Range("I122").Value = 4204.75
$ws_LTW.Range("J122").Value = 7494.5
$ws_LTW.Range("K122").Value = 12614.25
$ws_LTW.Range("L122").Value = 22483.5
$ws_LTW.Range("M122").Value = -10164.25
$ws_LTW.Range("N122").Value = -27383.5
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H19").Value = 2249
$ws_WVR.Range("J19").Value = 2249
$ws_WVR.Range("L19").Value = 2249
$ws_WVR.Range("N19").Value = -2597
$ws_WVR.Range("H42").Value = 49666
$ws_WVR.Range("I42").Value = 49999
$ws_WVR.Range("J42").Value = 49000
$ws_WVR.Range("K42").Value = 49999
$ws_WVR.Range("L42").Value = 49000
$ws_WVR.Range("M42").Value = -49621
$ws_WVR.Range("N42").Value = -49756
$ws_WVR.Range("H43").Value = 39999.5
$ws_WVR.Range("I43").Value = 49999
$ws_WVR.Range("J43").Value = 30000
$ws_WVR.Range("K43").Value = 49999
$ws_WVR.Range("L43").Value = 30000
$ws_WVR.Range("M43").Value = -49850
$ws_WVR.Range("N43").Value = -30298
$ws_WVR.Range("H103").Value = 0
$ws_WVR.Range("J103").Value = 0
$ws_WVR.Range("L103").Value = 0
$ws_WVR.Range("N103").ClearContents()
$ws_WVR.Range("H126").Value = 4583.2
$ws_WVR.Range("I126").Value = 2494.75
$ws_WVR.Range("K126").Value = 7484.25
$ws_WVR.Range("M126").Value = -5014.25
